# Settings back to qfin22 presentation - recalibration
$wb = $excel.ActiveWorkbook

# --- Update shared string array values on "Linear" worksheet ---
$wsLinear = $wb.Worksheets.Item("Linear")
$wsLinear.Range("B2").Value = -0.1343056305638665
$wsLinear.Range("B3").Value = 0.1405984146479144
$wsLinear.Range("B4").Value = 433.7702968771803
$wsLinear.Range("B5").Value = "[1.0, 0.25934787680264326, 0.08024270384055407, 0.1011318301121593, 0.060347240808209544, 0.06320652243483314, 0.18657242533447227, 0.29908203088843793, 0.17367868043586446, 0.02203711972053104, 0.034472153787815346, -0.005314975591231965, -0.01046895091000644, 0.13452919070520794, 0.2677602755667229, 0.11739406106603326, 0.0076051305339804, 0.02626120535685009, 0.02174826339204841, 0.002160365660303548]"

# --- Update values on "NonLinear" worksheet ---
$wsNonLinear = $wb.Worksheets.Item("NonLinear")
$wsNonLinear.Range("B4").Value = 0.01653149215147305
$wsNonLinear.Range("B5").Value = 0.1337848362156401
$wsNonLinear.Range("B6").Value = 452.8899395051749
$wsNonLinear.Range("B7").Value = -0.666474096740145
$wsNonLinear.Range("B8").Value = 0.1949508505165101
$wsNonLinear.Range("B9").Value = 416.1466970756001
$wsNonLinear.Range("B10").Value = "[1.0, 0.261205648326754, 0.08201540835261031, 0.10132809684962374, 0.0603468344942114, 0.06391436192950359, 0.18792716027141965, 0.3008810148951331, 0.17462184314773035, 0.022926409539673506, 0.034164155406926716, -0.004685031196414483, -0.009249855509447217, 0.135343447166886, 0.26806045507425347, 0.11820085764398873, 0.009209606610568253, 0.02590751932671112, 0.020470795560910433, 0.0030604113812347086]"
